$d = $word.ActiveDocument

$replacements = @(
    @("16×85=1360", "28×48=1344"),
    @("63×98=6174", "85×59=5015"),
    @("28×99=2772", "49×61=2989"),
    @("17×79=1343", "72×15=1080"),
    @("63×36=2268", "46×46=2116"),
    @("96×66=6336", "29×65=1885"),
    @("20×87=1740", "15×80=1200"),
    @("50×55=2750", "69×70=4830"),
    @("27×64=1728", "28×53=1484"),
    @("35×22=770",  "51×31=1581"),
    @("89×63=5607", "66×66=4356"),
    @("92×78=7176", "66×69=4554"),
    @("63×37=2331", "34×52=1768"),
    @("91×67=6097", "18×46=828"),
    @("27×31=837",  "42×78=3276"),
    @("71×30=2130", "71×52=3692"),
    @("45×95=4275", "18×90=1620"),
    @("37×28=1036", "22×38=836"),
    @("33×89=2937", "68×75=5100"),
    @("41×72=2952", "64×50=3200"),
    @("65×11=715",  "73×38=2774"),
    @("68×57=3876", "66×45=2970"),
    @("25×20=500",  "14×71=994"),
    @("18×78=1404", "18×73=1314"),
    @("77×68=5236", "93×41=3813")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
